$p = $ppt.ActivePresentation

# Slide 1 title: consolidate "First" + " " + "slide" runs into a single run.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = "placeholder_xyz"
$tr1.Text = "First slide"

# Slide 3 title: consolidate "Third" + " " + "slide" runs into a single run.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "placeholder_xyz"
$tr3.Text = "Third slide"
